# Add a new alias row mapping "CPSC 5115U" -> "CPSC 4115"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "CPSC 5115U"
$ws.Range("B6").Value = "CPSC 4115"

$ws.Range("B8").Select()
